# UV-5R Radio Configuration Table - add v1.2.0 sheet (PMR channels 1-16)
$wb = $excel.ActiveWorkbook

$ws11 = $wb.Worksheets.Item("v1.1.0")

# Add the new worksheet immediately after "v1.1.0" and rename it.
$ws12 = $wb.Worksheets.Add($null, $ws11)
$ws12.Name = "v1.2.0"

# Pre-format the used range the same way the other tables are formatted
# (center-aligned cells -> matches the existing shared cell style).
$ws12.Range("B2:T18").HorizontalAlignment = -4108

# ---- Header row ----
$ws12.Range("B2").Value = "Location"
$ws12.Range("C2").Value = "Name"
$ws12.Range("D2").Value = "Frequency"
$ws12.Range("E2").Value = "Duplex"
$ws12.Range("F2").Value = "Offset"
$ws12.Range("G2").Value = "Tone"
$ws12.Range("H2").Value = "rToneFreq"
$ws12.Range("I2").Value = "cToneFreq"
$ws12.Range("J2").Value = "DtcsCode"
$ws12.Range("K2").Value = "DtcsPolarity"
$ws12.Range("L2").Value = "Mode"
$ws12.Range("M2").Value = "TStep"
$ws12.Range("N2").Value = "Power"
$ws12.Range("O2").Value = "Skip"
$ws12.Range("P2").Value = "Comment"
$ws12.Range("Q2").Value = "URCALL"
$ws12.Range("R2").Value = "RPT1CALL"
$ws12.Range("S2").Value = "RPT2CALL"
$ws12.Range("T2").Value = "DVCODE"

# ---- Data rows: 16 PMR channels ----
$names = @("PMR 01","PMR 02","PMR 03","PMR 04","PMR 05","PMR 06","PMR 07","PMR 08","PMR 09","PMR 10","PMR 11","PMR 12","PMR 13","PMR 14","PMR 15","PMR 16")
$freqs = @(446.00625000000002,446.01875000000001,446.03125,446.04374999999999,446.05624999999998,446.06875000000002,446.08125000000001,446.09375,446.10624999999999,446.11874999999998,446.13125000000002,446.14375000000001,446.15625,446.16874999999999,446.18124999999998,446.19375000000002)

for ($i = 0; $i -lt 16; $i++) {
    $r = 3 + $i
    $ws12.Range("B$r").Value = ($i + 1)
    $ws12.Range("C$r").Value = $names[$i]
    $ws12.Range("D$r").Value = $freqs[$i]
    $ws12.Range("E$r").Value = $null
    $ws12.Range("F$r").Value = 0
    $ws12.Range("G$r").Value = $null
    $ws12.Range("H$r").Value = 88.5
    $ws12.Range("I$r").Value = 88.5
    $ws12.Range("J$r").Value = 23
    $ws12.Range("K$r").Value = "NN"
    $ws12.Range("L$r").Value = "NFM"
    $ws12.Range("M$r").Value = 5
    $ws12.Range("N$r").Value = "Low"
    $ws12.Range("O$r").Value = $null
    $ws12.Range("P$r").Value = $null
    $ws12.Range("Q$r").Value = $null
    $ws12.Range("R$r").Value = $null
    $ws12.Range("S$r").Value = $null
    $ws12.Range("T$r").Value = $null
}

# ---- Turn the range into a table, matching Table1 / Table2 ----
$lo = $ws12.ListObjects.Add(1, $ws12.Range("B2:T18"), $null, 1)

# Column widths (approximate best-fit sizing used on the other sheets).
$ws12.Columns.Item(2).ColumnWidth = 12
$ws12.Columns.Item(3).ColumnWidth = 9.91
$ws12.Columns.Item(4).ColumnWidth = 13.58
$ws12.Columns.Item(5).ColumnWidth = 10.75
$ws12.Columns.Item(6).ColumnWidth = 10.08
$ws12.Columns.Item(7).ColumnWidth = 9.08
$ws12.Columns.Item(8).ColumnWidth = 13.58
$ws12.Columns.Item(9).ColumnWidth = 13.66
$ws12.Columns.Item(10).ColumnWidth = 12.75
$ws12.Columns.Item(11).ColumnWidth = 15
$ws12.Columns.Item(12).ColumnWidth = 9.83
$ws12.Columns.Item(13).ColumnWidth = 9.58
$ws12.Columns.Item(14).ColumnWidth = 10.33
$ws12.Columns.Item(15).ColumnWidth = 8.41
$ws12.Columns.Item(16).ColumnWidth = 13.16
$ws12.Columns.Item(17).ColumnWidth = 11.41
$ws12.Columns.Item(18).ColumnWidth = 13.16
$ws12.Columns.Item(19).ColumnWidth = 13.16
$ws12.Columns.Item(20).ColumnWidth = 12

# Selection + make v1.2.0 the active tab, like the recorded session.
$ws12.Range("O8").Select()
$ws12.Activate()
